$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '60.187.88'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.31%  '
$ws.Range('E2').ClearFormats()

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.425.10'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.36%  '
$ws.Range('E3').ClearFormats()

$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('E4').ClearFormats()

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '553.94'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.47%  '
$ws.Range('E5').ClearFormats()

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '137.69'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -1.32%  '
$ws.Range('E6').ClearFormats()

$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E7').ClearFormats()

$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +1.83%  '
$ws.Range('E8').ClearFormats()

$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -0.88%  '
$ws.Range('E9').ClearFormats()

$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -1.10%  '
$ws.Range('E10').ClearFormats()

$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.18%  '
$ws.Range('E11').ClearFormats()

$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -1.21%  '
$ws.Range('E12').ClearFormats()

$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.45%  '
$ws.Range('E13').ClearFormats()

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.857.88'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -0.23%  '
$ws.Range('E14').ClearFormats()

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '60.107.04'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -0.31%  '
$ws.Range('E15').ClearFormats()

$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -1.14%  '
$ws.Range('E16').ClearFormats()

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.367.70'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -2.52%  '
$ws.Range('E17').ClearFormats()

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '11.34'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.68%  '
$ws.Range('E18').ClearFormats()

$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +1.96%  '
$ws.Range('E19').ClearFormats()

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '328.70'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -1.74%  '
$ws.Range('E20').ClearFormats()

$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.39%  '
$ws.Range('E21').ClearFormats()

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.00'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.01%  '
$ws.Range('E22').ClearFormats()

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '65.75'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.55%  '
$ws.Range('E23').ClearFormats()

$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +3.84%  '
$ws.Range('E24').ClearFormats()

$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +1.16%  '
$ws.Range('E25').ClearFormats()

$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.05%  '
$ws.Range('E26').ClearFormats()

$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +2.68%  '
$ws.Range('E27').ClearFormats()

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0₃0777'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -1.89%  '
$ws.Range('E28').ClearFormats()

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.77'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -1.13%  '
$ws.Range('E29').ClearFormats()

$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.83%  '
$ws.Range('E30').ClearFormats()

$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -3.22%  '
$ws.Range('E31').ClearFormats()

$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +3.68%  '
$ws.Range('E32').ClearFormats()

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.407'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -3.85%  '
$ws.Range('E33').ClearFormats()

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '18.63'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -0.87%  '
$ws.Range('E34').ClearFormats()

$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +1.34%  '
$ws.Range('E36').ClearFormats()

$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.32%  '
$ws.Range('E37').ClearFormats()

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '334.02'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +3.75%  '
$ws.Range('E39').ClearFormats()

$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -0.61%  '
$ws.Range('E40').ClearFormats()

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '38.87'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -2.52%  '
$ws.Range('E41').ClearFormats()

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '146.55'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +3.71%  '
$ws.Range('E42').ClearFormats()

$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -1.33%  '
$ws.Range('E43').ClearFormats()

$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +2.55%  '
$ws.Range('E44').ClearFormats()

$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.65%  '
$ws.Range('E45').ClearFormats()

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.580'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.90%  '
$ws.Range('E47').ClearFormats()

$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -1.56%  '
$ws.Range('E48').ClearFormats()

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '11.03'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.17%  '
$ws.Range('E49').ClearFormats()

$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -2.56%  '
$ws.Range('E50').ClearFormats()

$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -1.08%  '
$ws.Range('E51').ClearFormats()
